$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "C2" = 6.152135856756273
    "D2" = 5.280040698178554
    "E2" = 11.12338531991286
    "F2" = 49.44050186604019
    "G2" = 3.744214391672402
    "I2" = 35.29737795692802
    "J2" = 10.02788300527148
    "K2" = 19.59612556902693
    "M2" = 19.92064043291422
    "C3" = 6.158117380940759
    "D3" = 5.280372464145907
    "E3" = 11.14186304279861
    "F3" = 49.32812820919656
    "G3" = 3.748260643917604
    "I3" = 35.22646929845956
    "J3" = 10.05019365851183
    "K3" = 19.32433715630202
    "M3" = 19.83666893995177
    "C4" = 6.161869599419492
    "D4" = 5.280800247818137
    "E4" = 11.15454861024931
    "F4" = 49.27148954335699
    "G4" = 3.750872056787252
    "I4" = 35.19169160636365
    "J4" = 10.06495067215513
    "K4" = 19.16117186664459
    "M4" = 19.78955543299273
    "C5" = 6.163419009417644
    "D5" = 5.281031236331131
    "E5" = 11.16005535678564
    "F5" = 49.25151873955257
    "G5" = 3.751968289999107
    "I5" = 35.17972230524357
    "J5" = 10.07123052944826
    "K5" = 19.09569958051033
    "M5" = 19.77148861607417
    "C6" = 6.163677527170226
    "D6" = 5.281073023544898
    "E6" = 11.1609901263012
    "F6" = 49.24839043080848
    "G6" = 3.752152258827192
    "I6" = 35.17786781788624
    "J6" = 10.07228938101104
    "K6" = 19.08489194267818
    "M6" = 19.76855742008831
    "C7" = 6.161890412465518
    "D7" = 5.280803133137868
    "E7" = 11.15462151018409
    "F7" = 49.27120761815202
    "G7" = 3.750886710991348
    "I7" = 35.19152126660817
    "J7" = 10.0650342862396
    "K7" = 19.16028464819998
    "M7" = 19.78930717482421
    "C8" = 6.154181988599559
    "D8" = 5.280108770878323
    "E8" = 11.12947845644192
    "F8" = 49.39919248073249
    "G8" = 3.745583257528255
    "I8" = 35.27111103750194
    "J8" = 10.03535619230077
    "K8" = 19.50169363092119
    "M8" = 19.89077355956057
    "C9" = 6.139680977797187
    "D9" = 5.280511894228806
    "E9" = 11.09079450447163
    "F9" = 49.74800523989914
    "G9" = 3.736185087004823
    "I9" = 35.4966060355496
    "J9" = 9.985546430963826
    "K9" = 20.19685759657751
    "M9" = 20.12432008992398
    "C10" = 6.129379886718651
    "D10" = 5.281867270973177
    "E10" = 11.06883227760607
    "F10" = 50.06341333490582
    "G10" = 3.729882924950143
    "I10" = 35.70434295859377
    "J10" = 9.954055439275928
    "K10" = 20.71809957318556
    "M10" = 20.31595745773591
    "C11" = 6.124765516912751
    "D11" = 5.282710509236979
    "E11" = 11.06023987014963
    "F11" = 50.21956566383905
    "G11" = 3.727145039002432
    "I11" = 35.80788216083344
    "J11" = 9.940835720683197
    "K11" = 20.95642413405381
    "M11" = 20.40724840135637
    "C12" = 6.123028110053889
    "D12" = 5.283062136338924
    "E12" = 11.05718688484922
    "F12" = 50.28049840646567
    "G12" = 3.726126689111953
    "I12" = 35.84837734386645
    "J12" = 9.935988621026425
    "K12" = 21.04675675878964
    "M12" = 20.44238726916926
    "C13" = 6.123401854479646
    "D13" = 5.282984974555085
    "E13" = 11.05783547546439
    "F13" = 50.26729570522649
    "G13" = 3.726345191322415
    "I13" = 35.83959895940021
    "J13" = 9.937025463780769
    "K13" = 21.02729966469783
    "M13" = 20.434794550748
    "C14" = 6.124622381884267
    "D14" = 5.282738791846661
    "E14" = 11.05998467722848
    "F14" = 50.22454266847547
    "G14" = 3.727060890113793
    "I14" = 35.81118800764428
    "J14" = 9.94043376228829
    "K14" = 20.96385466912188
    "M14" = 20.41012802722952
    "C15" = 6.125371276272368
    "D15" = 5.282592197388563
    "E15" = 11.06132726259728
    "F15" = 50.19858910034208
    "G15" = 3.727501672560746
    "I15" = 35.79395269149457
    "J15" = 9.94254213753457
    "K15" = 20.92500113444318
    "M15" = 20.39509247767141
    "C16" = 6.129682858583116
    "D16" = 5.281816700128448
    "E16" = 11.06942192386483
    "F16" = 50.05346168140357
    "G16" = 3.730064437607368
    "I16" = 35.69775727965948
    "J16" = 9.954941621701142
    "K16" = 20.70254112415948
    "M16" = 20.31007225376024
    "C17" = 6.132345977220742
    "D17" = 5.281398809154061
    "E17" = 11.07474566928729
    "F17" = 49.96766232672357
    "G17" = 3.731669563948203
    "I17" = 35.64105123241433
    "J17" = 9.962831442112067
    "K17" = 20.56631436063657
    "M17" = 20.2589530346897
    "C18" = 6.133884497183317
    "D18" = 5.281179801140488
    "E18" = 11.07793938266979
    "F18" = 49.91950697335896
    "G18" = 3.732604939444322
    "I18" = 35.60928706645504
    "J18" = 9.967473535367979
    "K18" = 20.48807936221399
    "M18" = 20.22993927499093
    "C19" = 6.134406586052884
    "D19" = 5.281109324761363
    "E19" = 11.07904333666867
    "F19" = 49.90340808311759
    "G19" = 3.732923731880784
    "I19" = 35.59867888667949
    "J19" = 9.969063146127127
    "K19" = 20.46161354428645
    "M19" = 20.22018314151647
    "C20" = 6.132061786117802
    "D20" = 5.281441086608741
    "E20" = 11.07416532583591
    "F20" = 49.97667236447445
    "G20" = 3.731497438988714
    "I20" = 35.646999629704
    "J20" = 9.961980786231067
    "K20" = 20.58080426853882
    "M20" = 20.26435469862166
    "C21" = 6.124263615991879
    "D21" = 5.2828102269898
    "E21" = 11.05934795816627
    "F21" = 50.23705156852276
    "G21" = 3.726850172866022
    "I21" = 35.8194981714689
    "J21" = 9.939428349916952
    "K21" = 20.98248841718746
    "M21" = 20.41735792658927
    "C22" = 6.119224925118487
    "D22" = 5.28389324656123
    "E22" = 11.05083405601168
    "F22" = 50.41771159596963
    "G22" = 3.723920272556884
    "I22" = 35.93973100522076
    "J22" = 9.925615309199969
    "K22" = 21.24546368709256
    "M22" = 20.52065962002569
    "C23" = 6.121908989239764
    "D23" = 5.283298089549733
    "E23" = 11.05527112634456
    "F23" = 50.32033806864828
    "G23" = 3.725474232508356
    "I23" = 35.87487939317825
    "J23" = 9.932902858342608
    "K23" = 21.10509650212519
    "M23" = 20.46523085124276
    "C24" = 6.132190245660542
    "D24" = 5.281421906782883
    "E24" = 11.07442728475631
    "F24" = 49.97259527484273
    "G24" = 3.731575217530327
    "I24" = 35.64430775161425
    "J24" = 9.962365037176564
    "K24" = 20.57425311401653
    "M24" = 20.26191143548948
    "C25" = 6.143540332558203
    "D25" = 5.280215654827771
    "E25" = 11.10012405084373
    "F25" = 49.64320755329864
    "G25" = 3.738621121334428
    "I25" = 35.42818922280772
    "J25" = 9.99812423728955
    "K25" = 20.00658106764665
    "M25" = 20.05753924239749
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}